$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$v = $ws.Range("B543").Value
"B543=$v"
